# Generate Report for Handoff
# Update the "Latest Handoff/Handback Datetime" values for the
# 5633bf77-ee12-4584-b2b9-72824164816c file row (row 6) across all sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: column D holds "Latest Handoff Date" text
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D6").Value = "2016-30-20 06:30:58"

# zh-cn sheet: column E holds "Latest Handoff Datetime" text
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-20 06:30:55"

# de-de sheet: column E holds "Latest Handoff Datetime" text
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-20 06:30:58"
